# added ifo gdp component analysis preprocessing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revised historical values (restated GVA evaluation-series figures)
$ws.Range("B7").Value  = 87.29000000000001
$ws.Range("B23").Value = 89.25
$ws.Range("B30").Value = 93.59999999999999
$ws.Range("B40").Value = 96.59
$ws.Range("B44").Value = 98.98
$ws.Range("B45").Value = 99.38
$ws.Range("B48").Value = 101.22
$ws.Range("B52").Value = 103.67
$ws.Range("B53").Value = 104.2
$ws.Range("B55").Value = 104.35
$ws.Range("B56").Value = 104.65
$ws.Range("B59").Value = 104.77
$ws.Range("B68").Value = 105.91
$ws.Range("B69").Value = 105.76
$ws.Range("B72").Value = 106.02
$ws.Range("B73").Value = 105.27
$ws.Range("B75").Value = 106.4
$ws.Range("B76").Value = 105.59
$ws.Range("B77").Value = 105.44
$ws.Range("B80").Value = 105.29
$ws.Range("B81").Value = 105.19

# New row 82 for the latest data release (2025-08-15 / value 105.26)
# Copy the date-column formatting (border/font/number-format) from the row above
$ws.Range("A81").Copy()
$ws.Range("A82").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A82").Value = 45884
$ws.Range("B82").Value = 105.26
